$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create row 4 as a formatted copy of row 3 (same styles/row look) ---
$ws.Range("A3:I3").Copy($ws.Range("A4:I4"))

# G4 uses the plain (non-black-colored) style that F3/F4 use, not the
# style G3/H3 use -- copy just the number format/font/fill/border over.
$ws.Range("F3").Copy()
$ws.Range("G4").PasteSpecial(-4122)

# --- Fix up row 3 values ---
# I3: was an empty/zero numeric placeholder -> now text "A, 0"
$ws.Range("I3").Value = "A, 0"

# F3: trim the long allele list down to just "A3"
$ws.Range("F3").Value = "A3"

# G3: replace the multi-colour rich-text runs with plain text "A2"
$ws.Range("G3").Value = "A2"

# --- Fill in the new row 4 values ---
$ws.Range("A4").Value = "P22"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = "A2"
$ws.Range("D4").Value = "P12"
$ws.Range("E4").Value = "A"
$ws.Range("F4").Value = "A3"
$ws.Range("G4").Value = "A11"
# H4 keeps the same rich "B7, / DQ6, DQ5" text copied from H3
$ws.Range("I4").Value = 0

# --- Move the active selection from A4 to G4 ---
$ws.Range("G4").Select() | Out-Null
